$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = -21.92
$ws.Range("A14").Value = -21.742
$ws.Range("C15").Value = -13.182
$ws.Range("A16").Value = -22.163
$ws.Range("A21").Value = -20.432
$ws.Range("C21").Value = -12.617
$ws.Range("C22").Value = -12.5
$ws.Range("A23").Value = -20.246
$ws.Range("C24").Value = -12.198
$ws.Range("A25").Value = -21.839
$ws.Range("A26").Value = -22.129
$ws.Range("C27").Value = -13.431
$ws.Range("C28").Value = -12.985
$ws.Range("A29").Value = -21.265
$ws.Range("C36").Value = -13.045
$ws.Range("C39").Value = -12.82
$ws.Range("A40").Value = -20.016
$ws.Range("C45").Value = -13.122
$ws.Range("C48").Value = -11.205
$ws.Range("C49").Value = -13.482
$ws.Range("C52").Value = -11.655
$ws.Range("A53").Value = -21.836
$ws.Range("C53").Value = -12.789
$ws.Range("C54").Value = -12.91
$ws.Range("A57").Value = -22.129
$ws.Range("C57").Value = -13.734
$ws.Range("A59").Value = -22.523
$ws.Range("A65").Value = -21.484
$ws.Range("A69").Value = -21.519
$ws.Range("C70").Value = -11.983
$ws.Range("C71").Value = -11.362
$ws.Range("A79").Value = -21.008
$ws.Range("A83").Value = -21.963
$ws.Range("C86").Value = -13.883
$ws.Range("C87").Value = -13.329
$ws.Range("C89").Value = -13.376
$ws.Range("A91").Value = -20.744
$ws.Range("A93").Value = -21.508
$ws.Range("A100").Value = -22.276
$ws.Range("C101").Value = -12.721
$ws.Range("A103").Value = -21.902